$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPEC")

# Rename the PL2/SPPT header text (column E).
$ws.Cells.Item(1, 5).Value = "PL2/sPPT(W)"

# Insert a new column before column N (Resolution (Hz)) for the new "Display Type" field.
$ws.Columns.Item(14).Insert()

# New column header + values.
$ws.Cells.Item(1, 14).Value = "Display Type"
$ws.Cells.Item(2, 14).Value = "IPS"
$ws.Cells.Item(3, 14).Value = "IPS"

# vcore OC / vram OC values become "-" instead of 0.
$ws.Cells.Item(2, 10).Value = "-"
$ws.Cells.Item(3, 10).Value = "-"
$ws.Cells.Item(2, 11).Value = "-"
$ws.Cells.Item(3, 11).Value = "-"

$ws.Range("N1").Select()
